# Rakib & Mahfuz meal on from Noon, 12th March
#
# Sheet1 layout: row 2 holds day numbers (col B = day 1 ... col AF = day 31),
# so column M is day 12. Row 3 = Rakib, row 4 = Mahfuz (see column A / AK
# name labels). Both started taking the noon meal on the 12th, so their
# "day 12" meal count goes from 0 to 2. Every other changed cell in the
# workbook (AG/AH/AI/AJ per-person columns, the M17/AG18 totals, and the
# M44/AG44/M46/AG46 bazar-cost rows) is a formula that recalculates
# automatically from this single input change.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("M3").Value = 2
$ws.Range("M4").Value = 2

# Scroll the sheet view so row 18 / column F is the top-left visible cell,
# while leaving the current selection on M5 (matches the saved sheetView).
$excel.Goto($ws.Range("F18"), $false) | Out-Null
$win = $excel.ActiveWindow
$win.ScrollRow = 18
$win.ScrollColumn = 6
$ws.Range("M5").Select() | Out-Null
